$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert two new paragraphs at the very top of the document:
#      - Title style paragraph: "Resume - Ray Yan"
#      - Author style paragraph: "Ray Yan (Kin Long Yan)"
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "Resume - Ray Yan"
$titlePara.Style = "Title"

$secondPara = $d.Paragraphs.Item(2)
$secondPara.Range.InsertParagraphBefore()
$authorPara = $d.Paragraphs.Item(2)
$authorPara.Range.Text = "Ray Yan (Kin Long Yan)"
$authorPara.Style = "Author"

# ---------------------------------------------------------------------------
# 2. Add "Email:", "GitHub:", "LinkedIn:" labels before the corresponding
#    hyperlinks in the contact-info paragraph (now paragraph 3).
#    Insert from the right-most hyperlink back to the left-most one so that
#    earlier offsets stay valid while we work.
# ---------------------------------------------------------------------------
$hEmail = $d.Hyperlinks.Item(1)
$hGitHub = $d.Hyperlinks.Item(2)
$hLinkedIn = $d.Hyperlinks.Item(3)

$posEmail = $hEmail.Range.Start
$posGitHub = $hGitHub.Range.Start
$posLinkedIn = $hLinkedIn.Range.Start

$d.Range($posLinkedIn, $posLinkedIn).InsertBefore("LinkedIn: ")
$d.Range($posGitHub, $posGitHub).InsertBefore("GitHub: ")
$d.Range($posEmail, $posEmail).InsertBefore("Email: ")

# ---------------------------------------------------------------------------
# 3. Expand the CI/CD bullet point in the experience section with extra
#    detail about Jenkins Groovy scripting.
# ---------------------------------------------------------------------------
$oldText = "Design, implement, and maintain CI/CD pipelines using Jenkins and Bitbucket for deploying Web Applications, Android, and iOS applications in enterprise environments. Performed pipeline reconfiguration, ongoing maintenance, future pipeline development, and troubleshooting of deployment issues. Developed automation scripts using Shell, Python, and Jenkins Groovy scripts for CI/CD processes while implementing containerization solutions."
$newText = "Design, implement, and maintain CI/CD pipelines using Jenkins and Bitbucket for deploying Web Applications, Android, and iOS applications in enterprise environments. Developed complex Jenkins Groovy scripts for pipeline automation, custom build steps, and deployment orchestration. Performed pipeline reconfiguration, ongoing maintenance, future pipeline development, and troubleshooting of deployment issues. Created automation scripts using Shell, Python, and extensive Jenkins Groovy scripting for CI/CD processes while implementing containerization solutions."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

Write-Host "Done."
